# Auto-generated edit script for Linea 141 schedule update (commit: "Horarios actualizados Linea 141 - 259")
# Applies the new scraped rows / reordering described by the diff to all three worksheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$rows = @(
    @{r=2; A="Última actualización: 07:20:40"},
    @{r=3; A="Total filas: 84"},
    @{r=47; A="05:49:40"; B="07:04"; C="23_HERNANDEZ"; D=75; E="LP1912"},
    @{r=48; A="05:18:56"; B="07:04"; C="15_ABASTO"; D=106; E="LP1912"},
    @{r=57; A="07:20:40"; B="07:20"; C="10_OLMOS"; D=0; E="LP1912"},
    @{r=58; A="05:49:40"; B="07:21"; C="26_HERNANDEZ"; D=92; E="LP1912"},
    @{r=59; A="06:15:04"; B="07:23"; C="10_OLMOS"; D=68; E="LP1912"},
    @{r=60; A="05:49:40"; B="07:29"; C="10_OLMOS"; D=100; E="LP1912"},
    @{r=61; A="05:49:40"; B="07:31"; C="11_ETCHEVERRY"; D=102; E="LP1912"},
    @{r=62; A="05:49:40"; B="07:32"; C="84_COLONIA URQUIZA-ESC 49"; D=103; E="LP1912"},
    @{r=63; A="06:15:04"; B="07:32"; C="11_ETCHEVERRY"; D=77; E="LP1912"},
    @{r=64; A="07:20:40"; B="07:34"; C="23_HERNANDEZ"; D=14; E="LP1912"},
    @{r=65; A="05:49:40"; B="07:36"; C="27_EL RETIRO"; D=107; E="LP1912"},
    @{r=66; A="06:15:04"; B="07:37"; C="27_EL RETIRO"; D=82; E="LP1912"},
    @{r=67; A="05:49:40"; B="07:39"; C="10_OLMOS"; D=110; E="LP1912"},
    @{r=68; A="07:20:40"; B="07:46"; C="16_SANTA ANA"; D=26; E="LP1912"},
    @{r=69; A="06:43:40"; B="07:47"; C="14_ABASTO"; D=64; E="LP1912"},
    @{r=70; A="06:15:04"; B="07:48"; C="14_ABASTO"; D=93; E="LP1912"},
    @{r=71; A="06:43:40"; B="07:51"; C="215D_EL PATO"; D=68; E="LP1912"},
    @{r=72; A="06:15:04"; B="07:52"; C="215D_EL PATO"; D=97; E="LP1912"},
    @{r=73; A="07:20:40"; B="07:58"; C="16_SANTA ANA"; D=38; E="LP1912"},
    @{r=74; A="07:20:40"; B="07:59"; C="23_HERNANDEZ"; D=39; E="LP1912"},
    @{r=75; A="06:15:04"; B="08:01"; C="23_HERNANDEZ"; D=106; E="LP1912"},
    @{r=76; A="07:20:40"; B="08:03"; C="11_ETCHEVERRY"; D=43; E="LP1912"},
    @{r=77; A="06:43:40"; B="08:03"; C="23_HERNANDEZ"; D=80; E="LP1912"},
    @{r=78; A="06:57:30"; B="08:06"; C="23_HERNANDEZ"; D=69; E="LP1912"},
    @{r=79; A="06:15:04"; B="08:12"; C="15_ABASTO"; D=117; E="LP1912"},
    @{r=80; A="06:43:40"; B="08:21"; C="26_HERNANDEZ"; D=98; E="LP1912"},
    @{r=81; A="06:43:40"; B="08:22"; C="16_P MOR-SANTA ANA"; D=99; E="LP1912"},
    @{r=82; A="06:43:40"; B="08:23"; C="215B_EL PATO"; D=100; E="LP1912"},
    @{r=83; A="06:43:40"; B="08:27"; C="84_COLONIA URQUIZA-ESC 49"; D=104; E="LP1912"},
    @{r=84; A="06:43:40"; B="08:42"; C="81_EL PELIGRO"; D=119; E="LP1912"},
    @{r=85; A="07:20:40"; B="08:43"; C="14_ABASTO"; D=83; E="LP1912"},
    @{r=86; A="06:57:30"; B="08:54"; C="17_ROMERO"; D=117; E="LP1912"},
    @{r=87; A="07:20:40"; B="09:01"; C="215A_EL PATO"; D=101; E="LP1912"},
    @{r=88; A="07:20:40"; B="09:10"; C="16_P MOR-SANTA ANA"; D=110; E="LP1912"},
    @{r=89; A="07:20:40"; B="09:16"; C="27_EL RETIRO"; D=116; E="LP1912"}
)
foreach ($row in $rows) {
    if ($row.ContainsKey("A")) { $ws.Cells.Item($row.r, 1).Value = $row.A }
    if ($row.ContainsKey("B")) { $ws.Cells.Item($row.r, 2).Value = $row.B }
    if ($row.ContainsKey("C")) { $ws.Cells.Item($row.r, 3).Value = $row.C }
    if ($row.ContainsKey("D")) { $ws.Cells.Item($row.r, 4).Value = $row.D }
    if ($row.ContainsKey("E")) { $ws.Cells.Item($row.r, 5).Value = $row.E }
}

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$rows = @(
    @{r=2; A="Última actualización: 07:20:40"},
    @{r=3; A="Total filas: 14"},
    @{r=19; A="07:20:40"; B="09:01"; C="215A_EL PATO"; D=101; E="LP1912"}
)
foreach ($row in $rows) {
    if ($row.ContainsKey("A")) { $ws.Cells.Item($row.r, 1).Value = $row.A }
    if ($row.ContainsKey("B")) { $ws.Cells.Item($row.r, 2).Value = $row.B }
    if ($row.ContainsKey("C")) { $ws.Cells.Item($row.r, 3).Value = $row.C }
    if ($row.ContainsKey("D")) { $ws.Cells.Item($row.r, 4).Value = $row.D }
    if ($row.ContainsKey("E")) { $ws.Cells.Item($row.r, 5).Value = $row.E }
}

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$rows = @(
    @{r=2; A="Última actualización: 07:20:40"},
    @{r=3; A="Total filas: 19"},
    @{r=17; A="07:20:40"; B="07:37"; C="215A_LA PLATA"; D=17; E="L6173"},
    @{r=18; A="06:43:40"; B="08:06"; C="215C_LA PLATA"; D=83; E="L6203"},
    @{r=19; A="06:15:04"; B="08:07"; C="215C_LA PLATA"; D=112; E="L6203"},
    @{r=20; A="07:20:40"; B="08:09"; C="215C_LA PLATA"; D=49; E="L6203"},
    @{r=21; A="06:57:30"; B="08:10"; C="215C_LA PLATA"; D=73; E="L6203"},
    @{r=22; A="06:57:30"; B="08:35"; C="215A_LA PLATA"; D=98; E="L6173"},
    @{r=23; A="06:43:40"; B="08:38"; C="215A_LA PLATA"; D=115; E="L6173"},
    @{r=24; A="07:20:40"; B="09:08"; C="215D_LA PLATA"; D=108; E="L6203"}
)
foreach ($row in $rows) {
    if ($row.ContainsKey("A")) { $ws.Cells.Item($row.r, 1).Value = $row.A }
    if ($row.ContainsKey("B")) { $ws.Cells.Item($row.r, 2).Value = $row.B }
    if ($row.ContainsKey("C")) { $ws.Cells.Item($row.r, 3).Value = $row.C }
    if ($row.ContainsKey("D")) { $ws.Cells.Item($row.r, 4).Value = $row.D }
    if ($row.ContainsKey("E")) { $ws.Cells.Item($row.r, 5).Value = $row.E }
}
